$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (quantity / location fixes) ---
# Row 2: MAT_A / DC_001 - future production fix
$ws.Range("F2").Value = -233

# Row 3: MAT_A / DC_002 - future production fix
$ws.Range("F3").Value = -441

# Row 4: MAT_A / PLANT_001 - future production fix
$ws.Range("F4").Value = -706

# Row 5: was MAT_B / DC_002, becomes MAT_B / DC_001 with updated qty/horizon
$ws.Range("B5").Value = "DC_001"
$ws.Range("F5").Value = -103
$ws.Range("H5").Value = 4

# --- Append new rows 6 and 7 (previously missing DC_002 / PLANT_001 entries for MAT_B) ---

# Match the date number formatting used by the other rows' requirement_date / simulation_date columns
$ws.Range("C6").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("G6").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("C7").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("G7").NumberFormat = $ws.Range("G2").NumberFormat

# Row 6: MAT_B / DC_002
$ws.Range("A6").Value = "MAT_B"
$ws.Range("B6").Value = "DC_002"
$ws.Range("C6").Value = 45295
$ws.Range("D6").Value = "Distribution Demand - Forecast"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -33
$ws.Range("G6").Value = 45294
$ws.Range("H6").Value = 1

# Row 7: MAT_B / PLANT_001
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = 45295
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -103
$ws.Range("G7").Value = 45294
$ws.Range("H7").Value = 1
